$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 8800669077
$ws.Range("D4").Value = "Testing {#name#}"
$ws.Range("C4").Value = "text"
$ws.Range("B4").Value = "__BLANK__"

$ws.Range("E5").Select()
